$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A126").NumberFormat = "@"
$ws.Range("A126").Value = "01-07-2021"
$ws.Range("A126").Style = "Normal"
$ws.Range("B126").Value = 70473
$ws.Range("C126").Value = 10505

$ws.Range("A127").NumberFormat = "@"
$ws.Range("A127").Value = "02-07-2021"
$ws.Range("A127").Style = "Normal"
$ws.Range("B127").Value = 70825
$ws.Range("C127").Value = 10518

$ws.Range("A128").NumberFormat = "@"
$ws.Range("A128").Value = "05-07-2021"
$ws.Range("A128").Style = "Normal"
$ws.Range("B128").Value = 70599
$ws.Range("C128").Value = 10545

$ws.Range("A129").NumberFormat = "@"
$ws.Range("A129").Value = "06-07-2021"
$ws.Range("A129").Style = "Normal"
$ws.Range("B129").Value = 70595
$ws.Range("C129").Value = 10562

$ws.Range("A130").NumberFormat = "@"
$ws.Range("A130").Value = "07-07-2021"
$ws.Range("A130").Style = "Normal"
$ws.Range("B130").Value = 70683
$ws.Range("C130").Value = 10581

$ws.Range("A131").NumberFormat = "@"
$ws.Range("A131").Value = "08-07-2021"
$ws.Range("A131").Style = "Normal"
$ws.Range("B131").Value = 70433
$ws.Range("C131").Value = 10597

$ws.Range("A132").NumberFormat = "@"
$ws.Range("A132").Value = "09-07-2021"
$ws.Range("A132").Style = "Normal"
$ws.Range("B132").Value = 70553
$ws.Range("C132").Value = 10621

$ws.Range("A133").NumberFormat = "@"
$ws.Range("A133").Value = "12-07-2021"
$ws.Range("A133").Style = "Normal"
$ws.Range("B133").Value = 70241
$ws.Range("C133").Value = 10647

$ws.Range("A134").NumberFormat = "@"
$ws.Range("A134").Value = "13-07-2021"
$ws.Range("A134").Style = "Normal"
$ws.Range("B134").Value = 70421
$ws.Range("C134").Value = 10655

$ws.Range("A135").NumberFormat = "@"
$ws.Range("A135").Value = "14-07-2021"
$ws.Range("A135").Style = "Normal"
$ws.Range("B135").Value = 70265
$ws.Range("C135").Value = 10672

$ws.Range("A136").NumberFormat = "@"
$ws.Range("A136").Value = "15-07-2021"
$ws.Range("A136").Style = "Normal"
$ws.Range("B136").Value = 70402
$ws.Range("C136").Value = 10686

$ws.Range("A137").NumberFormat = "@"
$ws.Range("A137").Value = "19-07-2021"
$ws.Range("A137").Style = "Normal"
$ws.Range("B137").Value = 69982
$ws.Range("C137").Value = 10699

$ws.Range("A138").NumberFormat = "@"
$ws.Range("A138").Value = "20-07-2021"
$ws.Range("A138").Style = "Normal"
$ws.Range("B138").Value = 69682
$ws.Range("C138").Value = 10708

$ws.Range("A139").NumberFormat = "@"
$ws.Range("A139").Value = "21-07-2021"
$ws.Range("A139").Style = "Normal"
$ws.Range("B139").Value = 69740
$ws.Range("C139").Value = 10710

$ws.Range("A140").NumberFormat = "@"
$ws.Range("A140").Value = "22-07-2021"
$ws.Range("A140").Style = "Normal"
$ws.Range("B140").Value = 69826
$ws.Range("C140").Value = 10716

$ws.Range("A141").NumberFormat = "@"
$ws.Range("A141").Value = "23-07-2021"
$ws.Range("A141").Style = "Normal"
$ws.Range("B141").Value = 70031
$ws.Range("C141").Value = 10721

$ws.Range("A142").NumberFormat = "@"
$ws.Range("A142").Value = "26-07-2021"
$ws.Range("A142").Style = "Normal"
$ws.Range("B142").Value = 69642
$ws.Range("C142").Value = 10722

$ws.Range("A143").NumberFormat = "@"
$ws.Range("A143").Value = "27-07-2021"
$ws.Range("A143").Style = "Normal"
$ws.Range("B143").Value = 70008
$ws.Range("C143").Value = 10727

$ws.Range("A144").NumberFormat = "@"
$ws.Range("A144").Value = "28-07-2021"
$ws.Range("A144").Style = "Normal"
$ws.Range("B144").Value = 70156
$ws.Range("C144").Value = 10730

$ws.Range("A145").NumberFormat = "@"
$ws.Range("A145").Value = "29-07-2021"
$ws.Range("A145").Style = "Normal"
$ws.Range("B145").Value = 70208
$ws.Range("C145").Value = 10737

$ws.Range("A146").NumberFormat = "@"
$ws.Range("A146").Value = "30-07-2021"
$ws.Range("A146").Style = "Normal"
$ws.Range("B146").Value = 71852
$ws.Range("C146").Value = 10768
